$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Φύλλο1")

$ws.Range("A10").Value = "Polychronidou P"
$ws.Range("B10").Value = "IHU"
$ws.Range("C10").Value = "SOD"
$ws.Range("D10").Value = "Economic Sciences"
$ws.Range("E10").Value = "Assoc Professor"

$ws.Range("A10").Select()
